# "drugi problem drugega dne" - rename the existing sheet to "problem 1",
# add a new (second) sheet "problem 2" after it, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# Rename the original (only) sheet to "problem 1".
$ws1 = $wb.ActiveSheet
$ws1.Name = "problem 1"

# Add a brand-new worksheet right after "problem 1" and name it "problem 2".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "problem 2"

# Match the original sheet's default page-margin conventions.
$ws2.PageSetup.LeftMargin = 0.75 * 72
$ws2.PageSetup.RightMargin = 0.75 * 72
$ws2.PageSetup.TopMargin = 1 * 72
$ws2.PageSetup.BottomMargin = 1 * 72
$ws2.PageSetup.HeaderMargin = 0.5 * 72
$ws2.PageSetup.FooterMargin = 0.5 * 72

# The newly added sheet becomes the active / selected tab.
$ws2.Select()
